# Generate Report for Handback
#
# - Overview sheet: status "Ready for handoff" -> "Handed back: in sync with en-US"
#   for both language columns (zh-cn, de-de) on both data rows.
# - zh-cn / de-de sheets: fill in the "Latest Target File" (hyperlink to the
#   source .md on GitHub), "Latest Handback File" (the handed-back xlf, same
#   name as the handoff file since content is in sync) and "Latest Handback
#   DateTime" columns for both rows.
# - Column widths widen to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update the per-language status cells
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Hyperlink targets (same ones already used by column A on each sheet)
# ---------------------------------------------------------------------------
$file1Name = "30c1937e-af3c-4537-8f26-9b07f24af10f.md"
$file2Name = "a3f10386-b88f-4224-b418-e3ad96775d41.md"

$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77045efd0b02209e77186843a8991ea79a43c299/e2e/30c1937e-af3c-4537-8f26-9b07f24af10f.md"
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77045efd0b02209e77186843a8991ea79a43c299/e2e/a3f10386-b88f-4224-b418-e3ad96775d41.md"

# ---------------------------------------------------------------------------
# zh-cn sheet: Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = $file1Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $file1Url, "", "", $file1Name) | Out-Null
$wsZhCn.Range("J2").Value = "30c1937e-af3c-4537-8f26-9b07f24af10f.b8312d67dccb4b21ca5b28611daac098f84acd55.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-07 03:38:59"

$wsZhCn.Range("I3").Value = $file2Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $file2Url, "", "", $file2Name) | Out-Null
$wsZhCn.Range("J3").Value = "a3f10386-b88f-4224-b418-e3ad96775d41.bda0c7f29ed1e5037d9b3ff3f6bc8a436f75e85e.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-07 03:38:59"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $file1Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $file1Url, "", "", $file1Name) | Out-Null
$wsDeDe.Range("J2").Value = "30c1937e-af3c-4537-8f26-9b07f24af10f.b8312d67dccb4b21ca5b28611daac098f84acd55.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-07 03:39:21"

$wsDeDe.Range("I3").Value = $file2Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $file2Url, "", "", $file2Name) | Out-Null
$wsDeDe.Range("J3").Value = "a3f10386-b88f-4224-b418-e3ad96775d41.bda0c7f29ed1e5037d9b3ff3f6bc8a436f75e85e.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-07 03:39:21"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
